$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1319.7
$ws.Range("J12").Value = 1179.8
$ws.Range("L12").Value = 1179.8
$ws.Range("N12").Value = -1519.8
$ws.Range("H64").Value = 18571.428
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 18571.428
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H109").Value = 44899
$ws.Range("J109").Value = 44899
$ws.Range("L109").Value = 44899
$ws.Range("N109").Value = -47673
$ws.Range("H120").Value = 99499
$ws.Range("J120").Value = 99499
$ws.Range("L120").Value = 99499
$ws.Range("N120").Value = -109175
$ws.Range("H132").Value = 2321
$ws.Range("I132").Value = 2264.7778
$ws.Range("K132").Value = 6794.3334
$ws.Range("M132").Value = -4264.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3656.8333
$ws.Range("I2").Value = 481.33334
$ws.Range("J2").Value = 6832.3335
$ws.Range("K2").Value = 481.33334
$ws.Range("L2").Value = 6832.3335
$ws.Range("M2").Value = -368.33334
$ws.Range("N2").Value = -7058.3335
$ws.Range("H116").Value = 3656.8333
$ws.Range("I116").Value = 481.33334
$ws.Range("J116").Value = 6832.3335
$ws.Range("K116").Value = 481.33334
$ws.Range("L116").Value = 6832.3335
$ws.Range("M116").Value = 1812.66666
$ws.Range("N116").Value = -11420.3335
$ws.Range("H122").Value = 1770.3334
$ws.Range("I122").Value = 1155.5
$ws.Range("K122").Value = 3466.5
$ws.Range("M122").Value = -1016.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3656.8333
$ws.Range("I3").Value = 481.33334
$ws.Range("J3").Value = 6832.3335
$ws.Range("K3").Value = 481.33334
$ws.Range("L3").Value = 6832.3335
$ws.Range("M3").Value = -367.33334
$ws.Range("N3").Value = -7060.3335
$ws.Range("H86").Value = 8213.272000000001
$ws.Range("I86").Value = 6987
$ws.Range("K86").Value = 6987
$ws.Range("M86").Value = -5864
$ws.Range("H89").Value = 8213.272000000001
$ws.Range("I89").Value = 6987
$ws.Range("K89").Value = 34935
$ws.Range("M89").Value = -29319
$ws.Range("H134").Value = 2277.25
$ws.Range("I134").Value = 1171.4546
$ws.Range("J134").Value = 4710
$ws.Range("K134").Value = 3514.3638
$ws.Range("L134").Value = 14130
$ws.Range("M134").Value = -979.3638000000001
$ws.Range("N134").Value = -19200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 211
$ws.Range("I2").Value = 211
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 211
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -98
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 2589
$ws.Range("J7").Value = 5526.5
$ws.Range("L7").Value = 5526.5
$ws.Range("N7").Value = -5752.5
$ws.Range("H22").Value = 1472
$ws.Range("I22").Value = 660
$ws.Range("J22").Value = 1675
$ws.Range("K22").Value = 660
$ws.Range("L22").Value = 1675
$ws.Range("M22").Value = -310
$ws.Range("N22").Value = -2375
$ws.Range("H28").Value = 12745.2
$ws.Range("J28").Value = 12745.2
$ws.Range("L28").Value = 12745.2
$ws.Range("N28").Value = -13235.2
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H95").Value = 23998.334
$ws.Range("J95").Value = 23998.334
$ws.Range("L95").Value = 23998.334
$ws.Range("N95").Value = -29490.334
$ws.Range("H105").Value = 1547.875
$ws.Range("I105").Value = 1097
$ws.Range("J105").Value = 1998.75
$ws.Range("K105").Value = 1097
$ws.Range("L105").Value = 1998.75
$ws.Range("M105").Value = 650
$ws.Range("N105").Value = -5492.75
$ws.Range("H132").Value = 5736.125
$ws.Range("I132").Value = 5736.125
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 17208.375
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14678.375
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 222223090
$ws.Range("H80").Value = 9713.857
$ws.Range("I80").Value = 2666.6667
$ws.Range("J80").Value = 14999.25
$ws.Range("K80").Value = 8000.000100000001
$ws.Range("L80").Value = 44997.75
$ws.Range("M80").Value = -7064.000100000001
$ws.Range("N80").Value = -46869.75
$ws.Range("H83").Value = 9713.857
$ws.Range("I83").Value = 2666.6667
$ws.Range("J83").Value = 14999.25
$ws.Range("K83").Value = 24000.0003
$ws.Range("L83").Value = 134993.25
$ws.Range("M83").Value = -19320.0003
$ws.Range("N83").Value = -144353.25
$ws.Range("H107").Value = 896.8461
$ws.Range("I107").Value = 699
$ws.Range("K107").Value = 2097
$ws.Range("M107").Value = -177
$ws.Range("H131").Value = 435.6
$ws.Range("J131").Value = 595
$ws.Range("L131").Value = 1785
$ws.Range("N131").Value = -11865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9966.666999999999
$ws.Range("I70").Value = 9966.666999999999
$ws.Range("K70").Value = 9966.666999999999
$ws.Range("M70").Value = -9696.666999999999
$ws.Range("H73").Value = 9966.666999999999
$ws.Range("I73").Value = 9966.666999999999
$ws.Range("K73").Value = 9966.666999999999
$ws.Range("M73").Value = -9030.666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1280
$ws.Range("H22").Value = 1553.8462
$ws.Range("I22").Value = 1518.1818
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 1518.1818
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -1223.1818
$ws.Range("N22").Value = -2340
$ws.Range("H27").Value = 1553.8462
$ws.Range("I27").Value = 1518.1818
$ws.Range("J27").Value = 1750
$ws.Range("K27").Value = 1518.1818
$ws.Range("L27").Value = 1750
$ws.Range("M27").Value = -1411.1818
$ws.Range("N27").Value = -1964
$ws.Range("H43").Value = 16667
$ws.Range("J43").Value = 22499.5
$ws.Range("L43").Value = 22499.5
$ws.Range("N43").Value = -22885.5
$ws.Range("H122").Value = 6236.273
$ws.Range("I122").Value = 5199.875
$ws.Range("K122").Value = 15599.625
$ws.Range("M122").Value = -13149.625
$ws.Range("H132").Value = 5318.3
$ws.Range("I132").Value = 5454.7144
$ws.Range("K132").Value = 16364.1432
$ws.Range("M132").Value = -13834.1432
$ws.Range("H136").Value = 3272.818
$ws.Range("I136").Value = 3000.1
$ws.Range("K136").Value = 9000.299999999999
$ws.Range("M136").Value = -6450.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 32998.668
$ws.Range("J68").Value = 32998.668
$ws.Range("L68").Value = 32998.668
$ws.Range("N68").Value = -34620.668
$ws.Range("H71").Value = 32998.668
$ws.Range("J71").Value = 32998.668
$ws.Range("L71").Value = 98996.00399999999
$ws.Range("N71").Value = -107108.004
$ws.Range("H126").Value = 1750
$ws.Range("I126").Value = 1750
$ws.Range("K126").Value = 5250
$ws.Range("M126").Value = -2780
